$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "527.85", "0.999") are preserved as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "60.569.67"
$ws.Range("E2").Value = "  -1.82%  "

# Row 3
$ws.Range("D3").Value = "2.903.63"
$ws.Range("E3").Value = "  -2.84%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "527.85"
$ws.Range("E5").Value = "  -2.48%  "

# Row 6
$ws.Range("D6").Value = "142.71"
$ws.Range("E6").Value = "  -6.92%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").Value = "0.554"
$ws.Range("E8").Value = "  -1.77%  "

# Row 9
$ws.Range("D9").Value = "2.910.98"
$ws.Range("E9").Value = "  -2.72%  "

# Row 10
$ws.Range("E10").Value = "  -3.01%  "

# Row 11
$ws.Range("E11").Value = "  -4.39%  "

# Row 12
$ws.Range("D12").Value = "0.360"
$ws.Range("E12").Value = "  -0.98%  "

# Row 13
$ws.Range("D13").Value = "3.415.09"
$ws.Range("E13").Value = "  -2.79%  "

# Row 14
$ws.Range("E14").Value = "  +1.39%  "

# Row 15
$ws.Range("D15").Value = "60.573.75"
$ws.Range("E15").Value = "  -2.00%  "

# Row 16
$ws.Range("D16").Value = "22.60"
$ws.Range("E16").Value = "  -4.42%  "

# Row 17
$ws.Range("D17").Value = "2.907.14"
$ws.Range("E17").Value = "  -3.04%  "

# Row 18
$ws.Range("E18").Value = "  -3.43%  "

# Row 19
$ws.Range("D19").Value = "5.04"
$ws.Range("E19").Value = "  -1.10%  "

# Row 20
$ws.Range("D20").Value = "11.69"
$ws.Range("E20").Value = "  -1.77%  "

# Row 21
$ws.Range("D21").Value = "363.39"
$ws.Range("E21").Value = "  -6.73%  "

# Row 22
$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").Value = "64.15"
$ws.Range("E24").Value = "  -1.11%  "

# Row 25
$ws.Range("D25").Value = "3.023.87"
$ws.Range("E25").Value = "  -3.27%  "

# Row 26
$ws.Range("D26").Value = "0.452"
$ws.Range("E26").Value = "  -3.10%  "

# Row 27
$ws.Range("E27").Value = "  -3.53%  "

# Row 28
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  -7.13%  "

# Row 30
$ws.Range("D30").Value = [string]::Concat("0.0", [char]0x2083, "0861")
$ws.Range("E30").Value = "  -8.41%  "

# Row 31
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("E32").Value = "  -2.15%  "

# Row 33
$ws.Range("D33").Value = "19.60"
$ws.Range("E33").Value = "  -3.66%  "

# Row 34
$ws.Range("D34").Value = "148.05"
$ws.Range("E34").Value = "  -7.00%  "

# Row 35
$ws.Range("D35").Value = "4.36"
$ws.Range("E35").Value = "  -5.56%  "

# Row 36
$ws.Range("D36").Value = "5.59"
$ws.Range("E36").Value = "  -6.86%  "

# Row 37
$ws.Range("E37").Value = "  -5.62%  "

# Row 38
$ws.Range("E38").Value = "  -5.98%  "

# Row 39
$ws.Range("D39").Value = "37.93"
$ws.Range("E39").Value = "  +2.03%  "

# Row 40
$ws.Range("D40").Value = "1.50"
$ws.Range("E40").Value = "  -4.37%  "

# Row 41
$ws.Range("D41").Value = "2.328.85"
$ws.Range("E41").Value = "  -4.28%  "

# Row 42
$ws.Range("D42").Value = "3.67"
$ws.Range("E42").Value = "  -5.46%  "

# Row 43
$ws.Range("E43").Value = "  -2.23%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "20.65"
$ws.Range("E44").Value = "  -7.51%  "

# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0575"
$ws.Range("E45").Value = "  -3.23%  "

# Row 46
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").Value = "4.98"
$ws.Range("E47").Value = "  +1.33%  "

# Row 48
$ws.Range("D48").Value = "0.0235"
$ws.Range("E48").Value = "  -4.37%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.0933"
$ws.Range("E49").Value = "  -2.09%  "

# Row 50
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "10.36"
$ws.Range("E50").Value = "  -1.09%  "

# Row 51
$ws.Range("D51").Value = "18.47"
$ws.Range("E51").Value = "  -5.81%  "
